$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "23.420.21"
Set-TextValue $ws.Range("E2") "  +0.79%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.637.98"
Set-TextValue $ws.Range("E3") "  +2.22%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  -0.03%  "

# Row 5
Set-TextValue $ws.Range("E5") "  +0.04%  "

# Row 6
Set-TextValue $ws.Range("D6") "304.49"
Set-TextValue $ws.Range("E6") "  +0.31%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.3733"
Set-TextValue $ws.Range("E7") "  -1.17%  "

# Row 8
Set-TextValue $ws.Range("D8") "52.13"
Set-TextValue $ws.Range("E8") "  +0.69%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.3613"
Set-TextValue $ws.Range("E9") "  -0.56%  "

# Row 10
Set-TextValue $ws.Range("D10") "1.240"
Set-TextValue $ws.Range("E10") "  -2.73%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.08089"
Set-TextValue $ws.Range("E11") "  -0.58%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +0.09%  "

# Row 13
Set-TextValue $ws.Range("D13") "22.77"
Set-TextValue $ws.Range("E13") "  -0.35%  "

# Row 14
Set-TextValue $ws.Range("D14") "6.580"
Set-TextValue $ws.Range("E14") "  -0.44%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.00001265"
Set-TextValue $ws.Range("E15") "  +1.32%  "

# Row 16
Set-TextValue $ws.Range("D16") "7.265"
Set-TextValue $ws.Range("E16") "  -2.11%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.627.55"
Set-TextValue $ws.Range("E17") "  +1.29%  "

# Row 18
Set-TextValue $ws.Range("D18") "94.32"
Set-TextValue $ws.Range("E18") "  +0.33%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06877"
Set-TextValue $ws.Range("E19") "  -0.73%  "

# Row 20
Set-TextValue $ws.Range("D20") "18.07"
Set-TextValue $ws.Range("E20") "  -0.55%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.498"
Set-TextValue $ws.Range("E21") "  -0.61%  "

# Row 22
Set-TextValue $ws.Range("D22") "1.000"
Set-TextValue $ws.Range("E22") "  +0.00%  "

# Row 23
Set-TextValue $ws.Range("D23") "23.396.60"

# Row 24
Set-TextValue $ws.Range("D24") "12.72"
Set-TextValue $ws.Range("E24") "  -1.84%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.415"

# Row 26
Set-TextValue $ws.Range("D26") "3.006"
Set-TextValue $ws.Range("E26") "  -1.58%  "

# Row 27
Set-TextValue $ws.Range("D27") "21.13"
Set-TextValue $ws.Range("E27") "  -0.49%  "

# Row 28
Set-TextValue $ws.Range("D28") "151.55"
Set-TextValue $ws.Range("E28") "  +0.81%  "

# Row 29
Set-TextValue $ws.Range("D29") "5.332"
Set-TextValue $ws.Range("E29") "  +1.44%  "

# Row 30
Set-TextValue $ws.Range("D30") "135.11"
Set-TextValue $ws.Range("E30") "  +0.64%  "

# Row 31
Set-TextValue $ws.Range("D31") "2.273"
Set-TextValue $ws.Range("E31") "  -4.78%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.805.90"
Set-TextValue $ws.Range("E32") "  +1.25%  "

# Row 33
Set-TextValue $ws.Range("D33") "6.766"
Set-TextValue $ws.Range("E33") "  -0.13%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.9466"
Set-TextValue $ws.Range("E34") "  -2.11%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.02821"
Set-TextValue $ws.Range("E35") "  +2.68%  "

# Row 36
Set-TextValue $ws.Range("D36") "10.29"
Set-TextValue $ws.Range("E36") "  +0.26%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.2521"
Set-TextValue $ws.Range("E37") "  -0.79%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.07181"
Set-TextValue $ws.Range("E38") "  -4.71%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.08767"
Set-TextValue $ws.Range("E39") "  -0.64%  "

# Row 40
Set-TextValue $ws.Range("D40") "6.040"
Set-TextValue $ws.Range("E40") "  -1.39%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.364"
Set-TextValue $ws.Range("E41") "  -1.90%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.7007"
Set-TextValue $ws.Range("E42") "  -1.70%  "

# Row 43
Set-TextValue $ws.Range("D43") "12.39"
Set-TextValue $ws.Range("E43") "  -1.23%  "

# Row 44
Set-TextValue $ws.Range("D44") "15.97"
Set-TextValue $ws.Range("E44") "  +2.67%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.6480"
Set-TextValue $ws.Range("E45") "  -1.01%  "

# Row 46
Set-TextValue $ws.Range("B46") "NEARProtocol"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D46") "2.316"
Set-TextValue $ws.Range("E46") "  -0.42%  "

# Row 47
Set-TextValue $ws.Range("B47") "Frax"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D47") "1.000"
Set-TextValue $ws.Range("E47") "  +0.07%  "

# Row 48
Set-TextValue $ws.Range("D48") "3.999"

# Row 49
Set-TextValue $ws.Range("D49") "0.07965"
Set-TextValue $ws.Range("E49") "  +0.10%  "

# Row 50
Set-TextValue $ws.Range("D50") "128.21"
Set-TextValue $ws.Range("E50") "  -3.48%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.193"
Set-TextValue $ws.Range("E51") "  -1.03%  "
